$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.763.77'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.31%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.040.50'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.08%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.49'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.12%  '

$ws.Range('E6').Value = '  -0.93%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.69'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.11%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').Value = '  -2.89%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0842'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.15%  '

$ws.Range('E11').Value = '  -0.31%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.340.96'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.08%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.41'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.80%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.05'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.06%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.47'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +5.30%  '

$ws.Range('E16').Value = '  +0.61%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.046.84'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.34%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.712.01'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.36%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.42'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.65%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.91'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.77%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0826'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.15%  '

$ws.Range('E22').Value = '  -0.81%  '

$ws.Range('E23').Value = '  +0.56%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.43'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.84%  '

$ws.Range('E25').Value = '  +2.89%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.22'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.41%  '

$ws.Range('E27').Value = '  +0.91%  '

$ws.Range('E28').Value = '  -0.46%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.78'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.91%  '

$ws.Range('E30').Value = '  -0.05%  '

$ws.Range('E31').Value = '  -1.15%  '

$ws.Range('E32').Value = '  +8.85%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.38'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.65%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.49'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.31%  '

$ws.Range('E35').Value = '  -0.22%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.54'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.44%  '

$ws.Range('E37').Value = '  +3.75%  '

$ws.Range('E38').Value = '  +6.11%  '

$ws.Range('E39').Value = '  -0.20%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.99'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.37%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.527.89'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.80%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.62'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.51%  '

$ws.Range('E43').Value = '  -1.19%  '

$ws.Range('E44').Value = '  +0.01%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.20'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.91%  '

$ws.Range('E47').Value = '  -0.32%  '

$ws.Range('E48').Value = '  -0.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.10'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.41%  '

$ws.Range('E50').Value = '  -0.80%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.230.89'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.04%  '
